$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.967.09'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.949.10'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'378.53"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").Value = "'101.53"
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = "'0.0849"
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '3.414.86'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'7.77"
$ws.Range("E14").Value = '  +4.99%  '
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = "'12.19"
$ws.Range("E16").Value = '  +70.67%  '
$ws.Range("D17").Value = '2.948.32'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = '50.926.10'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  -3.71%  '
$ws.Range("D21").Value = "'12.37"
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = "'69.39"
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("D24").Value = "'266.27"
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = "'3.19"
$ws.Range("E25").Value = '  +8.79%  '
$ws.Range("E26").Value = '  -3.11%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").Value = "'7.04"
$ws.Range("E28").Value = '  -8.56%  '
$ws.Range("D29").Value = "'25.61"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").Value = "'0.164"
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("E31").Value = '  -4.48%  '
$ws.Range("D32").Value = "'10.19"
$ws.Range("E32").Value = '  +3.72%  '
$ws.Range("D33").Value = "'50.50"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = "'2.05"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = "'33.66"
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D36").Value = "'0.0431"
$ws.Range("E36").Value = '  -5.32%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = '  +4.40%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").Value = "'16.55"
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").Value = "'2.52"
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").Value = "'118.42"
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").Value = "'3.55"
$ws.Range("E44").Value = '  +9.36%  '
$ws.Range("D45").Value = "'21.35"
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").Value = '2.001.87'
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("D49").Value = "'0.258"
$ws.Range("E49").Value = '  -4.72%  '
$ws.Range("E50").Value = '  -9.35%  '
$ws.Range("D51").Value = "'5.29"
$ws.Range("E51").Value = '  +4.03%  '
